# Updated cryptos list on Mon May 29 13:28:17 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "28.012.60"
Set-TextValue $ws.Range("E2") "  +2.15%  "
Set-TextValue $ws.Range("D3") "1.907.22"
Set-TextValue $ws.Range("E3") "  +2.54%  "
Set-TextValue $ws.Range("D4") "1.004"
Set-TextValue $ws.Range("E4") "  -0.70%  "
Set-TextValue $ws.Range("D5") "317.78"
Set-TextValue $ws.Range("E5") "  +2.02%  "
Set-TextValue $ws.Range("E6") "  -0.76%  "
Set-TextValue $ws.Range("D7") "0.4822"
Set-TextValue $ws.Range("E7") "  +1.13%  "
Set-TextValue $ws.Range("D8") "0.3805"
Set-TextValue $ws.Range("E8") "  -0.14%  "
Set-TextValue $ws.Range("D9") "0.07366"
Set-TextValue $ws.Range("E9") "  +0.77%  "
Set-TextValue $ws.Range("E10") "  +0.26%  "
Set-TextValue $ws.Range("D11") "20.82"
Set-TextValue $ws.Range("E11") "  +0.12%  "
Set-TextValue $ws.Range("E12") "  -0.47%  "
Set-TextValue $ws.Range("D13") "1.907.28"
Set-TextValue $ws.Range("E13") "  +2.47%  "
Set-TextValue $ws.Range("D14") "5.492"
Set-TextValue $ws.Range("E14") "  +1.00%  "
Set-TextValue $ws.Range("D15") "6.647"
Set-TextValue $ws.Range("E15") "  +1.67%  "
Set-TextValue $ws.Range("D16") "91.80"
Set-TextValue $ws.Range("E16") "  +1.91%  "
Set-TextValue $ws.Range("D17") "1.004"
Set-TextValue $ws.Range("E17") "  -0.76%  "
Set-TextValue $ws.Range("D18") "0.000008896"
Set-TextValue $ws.Range("E18") "  +0.97%  "
Set-TextValue $ws.Range("E19") "  -0.67%  "
Set-TextValue $ws.Range("D20") "28.039.85"
Set-TextValue $ws.Range("E20") "  +2.37%  "
Set-TextValue $ws.Range("D21") "14.72"
Set-TextValue $ws.Range("E21") "  +0.60%  "
Set-TextValue $ws.Range("D22") "5.144"
Set-TextValue $ws.Range("D23") "2.156.82"
Set-TextValue $ws.Range("E23") "  +3.95%  "
Set-TextValue $ws.Range("D24") "10.90"
Set-TextValue $ws.Range("E24") "  +2.07%  "
Set-TextValue $ws.Range("D25") "156.70"
Set-TextValue $ws.Range("E25") "  +0.72%  "
Set-TextValue $ws.Range("E26") "  -1.73%  "
Set-TextValue $ws.Range("D27") "18.51"
Set-TextValue $ws.Range("E27") "  +0.30%  "
Set-TextValue $ws.Range("D28") "2.116"
Set-TextValue $ws.Range("E28") "  +5.44%  "
Set-TextValue $ws.Range("D29") "117.63"
Set-TextValue $ws.Range("E29") "  +2.09%  "
Set-TextValue $ws.Range("D30") "4.977"
Set-TextValue $ws.Range("E30") "  +0.52%  "
Set-TextValue $ws.Range("D31") "0.08946"
Set-TextValue $ws.Range("E31") "  +0.57%  "
Set-TextValue $ws.Range("D32") "3.239"
Set-TextValue $ws.Range("E32") "  -2.49%  "
Set-TextValue $ws.Range("D33") "1.255"
Set-TextValue $ws.Range("E33") "  +4.16%  "
Set-TextValue $ws.Range("D34") "0.7723"
Set-TextValue $ws.Range("E34") "  +2.91%  "
Set-TextValue $ws.Range("D35") "4.667"
Set-TextValue $ws.Range("E35") "  +1.66%  "
Set-TextValue $ws.Range("D36") "2.613"
Set-TextValue $ws.Range("E36") "  -3.93%  "
Set-TextValue $ws.Range("D37") "0.02049"
Set-TextValue $ws.Range("E37") "  +0.28%  "
Set-TextValue $ws.Range("D38") "1.105"
Set-TextValue $ws.Range("E38") "  -1.47%  "
Set-TextValue $ws.Range("E39") "  -0.40%  "
Set-TextValue $ws.Range("D40") "0.05290"
Set-TextValue $ws.Range("E40") "  +0.54%  "
Set-TextValue $ws.Range("D41") "2.996"
Set-TextValue $ws.Range("E41") "  +0.30%  "
Set-TextValue $ws.Range("D42") "6.995"
Set-TextValue $ws.Range("E42") "  -0.44%  "
Set-TextValue $ws.Range("B43") "Aptos"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D43") "8.507"
Set-TextValue $ws.Range("E43") "  -1.13%  "
Set-TextValue $ws.Range("B44") "Algorand"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D44") "0.1526"
Set-TextValue $ws.Range("E44") "  +0.21%  "
Set-TextValue $ws.Range("D45") "110.30"
Set-TextValue $ws.Range("E45") "  +7.38%  "
Set-TextValue $ws.Range("D46") "10.66"
Set-TextValue $ws.Range("E46") "  -0.29%  "
Set-TextValue $ws.Range("D47") "0.4831"
Set-TextValue $ws.Range("E47") "  -0.79%  "
Set-TextValue $ws.Range("D48") "1.003"
Set-TextValue $ws.Range("E48") "  -0.77%  "
Set-TextValue $ws.Range("D49") "1.646"
Set-TextValue $ws.Range("E49") "  -1.15%  "
Set-TextValue $ws.Range("D50") "68.05"
Set-TextValue $ws.Range("E50") "  +1.16%  "
Set-TextValue $ws.Range("D51") "0.06079"
Set-TextValue $ws.Range("E51") "  -0.23%  "
